$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.445.72"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "3.078.62"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.88"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.05"
$ws.Range("E6").Value = "  +5.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.075.32"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.12"
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("E12").Value = "  +4.60%  "
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.25"
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("D15").Value = "3.572.70"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "64.505.20"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").Value = "3.078.01"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.25"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.94"
$ws.Range("E21").Value = "  +2.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.57"
$ws.Range("E23").Value = "  +6.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.71"
$ws.Range("E24").Value = "  +10.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.80"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.81"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.15"
$ws.Range("E28").Value = "  +4.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.08"
$ws.Range("E29").Value = "  +5.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.26"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.50"
$ws.Range("E33").Value = "  +4.19%  "
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.20"
$ws.Range("E35").Value = "  +4.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.85"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "461.45"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  +19.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0831"
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0406"
$ws.Range("E40").Value = "  +3.65%  "
$ws.Range("D41").Value = "2.972.45"
$ws.Range("E41").Value = "  -4.42%  "
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.115"
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.82"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.264"
$ws.Range("E45").Value = "  +5.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  +6.21%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("E48").Value = "  +2.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "120.51"
$ws.Range("E49").Value = "  +3.89%  "
$ws.Range("D50").Value = "0.0₃0518"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("E51").Value = "  +1.60%  "
